$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 126: new problem entry "Ways to Express an Integer as Sum of Powers" ---
$ws.Range("A126").Value = 2787
$ws.Range("B126").Value = "Ways to Express an Integer as Sum of Powers"
$ws.Range("C126").Value = "#dynamic-programming "
$ws.Range("D126").Value = "medium"
$ws.Range("E126").Value = 0
$ws.Range("F126").Value = 1
$ws.Range("G126").Value = 30
$ws.Range("H126").Value = 45881
$ws.Range("I126").Value = 45881
$ws.Range("J126").Value = "?"
$ws.Rows.Item(126).RowHeight = 34

# --- Row 127: new problem entry "Power of Three" ---
$ws.Range("A127").Value = 326
$ws.Range("B127").Value = "Power of Three"
$ws.Range("C127").Value = "#math"
$ws.Range("D127").Value = "easy"
$ws.Range("E127").Value = 1
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 5
$ws.Range("H127").Value = 45882
$ws.Range("I127").Value = 45882
$ws.Rows.Item(127).RowHeight = 17

# --- Rows 128-129: trailing date-only rows (notes placeholders) ---
$ws.Range("H128").Value = 45883
$ws.Range("I128").Value = 45883
$ws.Range("H129").Value = 45883
$ws.Range("I129").Value = 45883

# The H/I (date) columns default to the plain centered style; copy the
# existing date-formatted style (m/d/yyyy, centered) from H2 down onto the
# new rows so the new cells share the same style index as the rest of the
# date column.
$ws.Range("H2").Copy()
$ws.Range("H126:I129").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to mirror where the user ended up after typing the
# last row of notes.
$null = $ws.Range("H129:I129").Select()
